# Daily attendance processing - 2025-12-24 11:54:24
# Normalizes the ordering of names/emails in the "Recorded By" column (G)
# so that automated/system accounts are listed in a consistent order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Exact, known "old value" -> "new value" replacements for column G
# (derived from the authoritative edit being applied to this report).
$map = @{
    "backup@backdoor.com, system, System" = "backup@backdoor.com, System, system";
    "dnasr281@gmail.com, System"          = "System, dnasr281@gmail.com";
    "System, backup@backdoor.com"         = "backup@backdoor.com, System";
    "dnasr281@gmail.com, admin@admin.com" = "admin@admin.com, dnasr281@gmail.com";
}

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $current = $cell.Value2
    if ($map.ContainsKey($current)) {
        $cell.Value2 = $map[$current]
    }
}
